# "Hover text on image functionality added"
#
# Adds a "NULL" placeholder hover-text value (column I / "DIS1") for every
# row that has a picture in column F but no hover text yet, and updates the
# sheet's current selection to the last row that was touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsNeedingHoverText = @(2, 5, 6, 8, 23, 29, 31, 32, 33)
foreach ($r in $rowsNeedingHoverText) {
    $ws.Cells.Item($r, 9).Value = "NULL"
}

# Scroll/selection bookkeeping to match where the author ended up working.
$ws.Activate() | Out-Null
$ws.Range("A4").Select() | Out-Null
$ws.Range("I32:I33").Select() | Out-Null
